$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A90").Value = "15-12-2025"
$ws.Range("B90").Value = "The price of gold in India today is ₹13,473 per gram for 24 karat gold, ₹12,350 per gram for 22 karat gold and ₹10,105 per gram for 18 karat gold (also called 999 gold)."
